$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.821.72'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '2.675.46'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.52'
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.71'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +3.62%  '
$ws.Range("E9").Value = '  +2.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.401'
$ws.Range("E10").Value = '  -0.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.87'
$ws.Range("E11").Value = '  -2.89%  '
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("E13").Value = '  -4.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.15'
$ws.Range("E14").Value = '  -3.48%  '
$ws.Range("D15").Value = '3.156.17'
$ws.Range("E15").Value = '  -0.61%  '
$ws.Range("D16").Value = '65.670.41'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").Value = '2.666.84'
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.77'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("E19").Value = '  -1.71%  '
$ws.Range("E20").Value = '  -3.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.93'
$ws.Range("E21").Value = '  -1.70%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.57'
$ws.Range("E23").Value = '  -2.52%  '
$ws.Range("E24").Value = '  +3.21%  '
$ws.Range("E25").Value = '  -1.60%  '
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.61'
$ws.Range("E27").Value = '  -3.51%  '
$ws.Range("E28").Value = '  -3.74%  '
$ws.Range("E29").Value = '  -3.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  -3.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '533.68'
$ws.Range("E32").Value = '  -1.37%  '
$ws.Range("E33").Value = '  -1.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.52'
$ws.Range("E35").Value = '  -0.41%  '
$ws.Range("E36").Value = '  -2.89%  '
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '158.14'
$ws.Range("E39").Value = '  -4.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.94'
$ws.Range("E40").Value = '  -2.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '163.53'
$ws.Range("E42").Value = '  -3.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.11'
$ws.Range("E43").Value = '  -1.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.39'
$ws.Range("E44").Value = '  +2.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0616'
$ws.Range("E45").Value = '  -2.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.81'
$ws.Range("E46").Value = '  -4.77%  '
$ws.Range("E47").Value = '  -3.39%  '
$ws.Range("D49").Value = '0.0₆0256'
$ws.Range("E49").Value = '  +8.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.13'
$ws.Range("E50").Value = '  -4.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0991'
$ws.Range("E51").Value = '  -0.50%  '
